$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows at position 4, pushing the old header row (5) and
#    data row (6) down to 7 and 8. Excel auto-adjusts existing formulas
#    (e.g. the A2/B2.. "ISNUMBER(A6)" refs become "ISNUMBER(A8)").
$ws.Rows("4:5").Insert()

# 2. Rename three (groups of) headers in what is now row 7.
$ws.Range("C7").Value = "LOSS20HALF"
$ws.Range("E7").Value = "VORSCHAEDEN_ANZAHL"
$ws.Range("AK7").Value = "VORSCHAEDEN0_typeKH"
$ws.Range("AL7").Value = "VORSCHAEDEN0_typetk"
$ws.Range("AM7").Value = "VORSCHAEDEN0_month"
$ws.Range("AN7").Value = "VORSCHAEDEN0_year"
$ws.Range("AO7").Value = "VORSCHAEDEN1_typetk"
$ws.Range("AP7").Value = "VORSCHAEDEN1_month"
$ws.Range("AQ7").Value = "VORSCHAEDEN1_year"
$ws.Range("AR7").Value = "VORSCHAEDEN2_typevk"
$ws.Range("AS7").Value = "VORSCHAEDEN2_month"
$ws.Range("AT7").Value = "VORSCHAEDEN2_year"

# 3. Small data fix: kofferraumvolumen_num goes from an integer to a float.
$ws.Range("BZ8").Value = 45.1

# 4. New row 4: quoted, comma-suffixed header names, read off row 7,
#    mirroring the existing row2/row3 "shared formula down columns" idiom.
$ws.Range("B4:CK4").Formula = '=+""""&B7&""""&","'

# 5. New row 5: the 1-based column number of each column.
$ws.Range("A5:CK5").Formula = "=+COLUMN()"
